# Update cryptos list prices (column D) and volume/1h percentages (column E)
# to reflect the latest scraped snapshot. Values that look like plain
# decimal numbers are prefixed with an apostrophe so Excel stores them
# as text (preserving trailing zeros / exact formatting) instead of
# silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "29.384.39"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.883.31"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'0.7122"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'242.37"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.08036"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'25.30"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("D12").Value = "1.892.60"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'0.7206"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "'5.245"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'92.98"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "'6.349"
$ws.Range("E16").Value = "  +5.71%  "
$ws.Range("D17").Value = "'0.000008459"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "29.399.39"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "2.151.45"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'241.41"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.849"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'0.1585"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'164.03"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'9.049"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'1.508"
$ws.Range("D30").Value = "'4.417"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'4.342"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "'1.197"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("D33").Value = "'0.05369"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'1.950"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "1.292.25"
$ws.Range("E38").Value = "  +10.02%  "
$ws.Range("D39").Value = "'0.01888"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").Value = "'2.747"
$ws.Range("D41").Value = "'6.600"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("D42").Value = "'0.9207"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "'111.97"
$ws.Range("E43").Value = "  +5.35%  "
$ws.Range("D44").Value = "'74.09"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("E46").Value = "  +5.62%  "
$ws.Range("D47").Value = "2.047.57"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").Value = "'1.808"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "'0.5220"
$ws.Range("D50").Value = "'9.499"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'0.4398"
$ws.Range("E51").Value = "  +2.07%  "